$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in column B for rows 2-5
$ws.Range("B2").Value = 444
$ws.Range("B3").Value = 445
$ws.Range("B4").Value = 446
$ws.Range("B5").Value = 447

# Remove row 6 entirely (was A6=1, B6=404)
$ws.Range("A6:B6").Delete()
